$wb = $excel.ActiveWorkbook

# Remember which sheet is currently active so we can restore it after
# touching the "By-Author" sheet's selection below.
$originalActiveSheet = $wb.ActiveSheet

# Rename sheets: "Per-Month" -> "By-Month", "Per-Author" -> "By-Author"
$wb.Worksheets.Item("Per-Month").Name = "By-Month"
$wb.Worksheets.Item("Per-Author").Name = "By-Author"

# Reset the lingering cell selection (previously E17) on the renamed
# "By-Author" sheet back to A1.
$wsAuthor = $wb.Worksheets.Item("By-Author")
$wsAuthor.Activate()
$wsAuthor.Range("A1").Select()

# Restore the originally active sheet/tab.
$originalActiveSheet.Activate()
